$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price/volume data per latest scrape
# Row 2
$ws.Range("D2").Value = "'76.471.06"
$ws.Range("E2").Value = "  +0.37%  "

# Row 3
$ws.Range("D3").Value = "'3.051.06"
$ws.Range("E3").Value = "  +3.49%  "

# Row 4
$ws.Range("E4").Value = "  +0.01%  "

# Row 5
$ws.Range("D5").Value = "'200.62"
$ws.Range("E5").Value = "  -1.95%  "

# Row 6
$ws.Range("D6").Value = "'624.69"
$ws.Range("E6").Value = "  +4.28%  "

# Row 7
$ws.Range("E7").Value = "  +0.05%  "

# Row 8
$ws.Range("E8").Value = "  -0.59%  "

# Row 9
$ws.Range("E9").Value = "  +3.50%  "

# Row 10
$ws.Range("D10").Value = "'3.050.32"
$ws.Range("E10").Value = "  +3.59%  "

# Row 11
$ws.Range("E11").Value = "  +0.40%  "

# Row 12
$ws.Range("E12").Value = "  -0.71%  "

# Row 13
$ws.Range("D13").Value = "'5.24"
$ws.Range("E13").Value = "  +5.75%  "

# Row 14
$ws.Range("D14").Value = "'3.611.64"
$ws.Range("E14").Value = "  +3.76%  "

# Row 15
$ws.Range("D15").Value = "'29.04"
$ws.Range("E15").Value = "  +2.67%  "

# Row 16
$ws.Range("D16").Value = "'76.366.66"
$ws.Range("E16").Value = "  +0.40%  "

# Row 17
$ws.Range("E17").Value = "  +1.56%  "

# Row 18
$ws.Range("D18").Value = "'3.044.94"
$ws.Range("E18").Value = "  +3.64%  "

# Row 19
$ws.Range("D19").Value = "'13.57"
$ws.Range("E19").Value = "  +2.17%  "

# Row 20
$ws.Range("D20").Value = "'9.05"
$ws.Range("E20").Value = "  +1.75%  "

# Row 21
$ws.Range("D21").Value = "'377.88"
$ws.Range("E21").Value = "  +0.99%  "

# Row 22
$ws.Range("D22").Value = "'2.32"
$ws.Range("E22").Value = "  +0.10%  "

# Row 23
$ws.Range("D23").Value = "'4.38"
$ws.Range("E23").Value = "  +1.63%  "

# Row 24
$ws.Range("D24").Value = "'73.36"
$ws.Range("E24").Value = "  +2.06%  "

# Row 26
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  +0.00%  "

# Row 27
$ws.Range("E27").Value = "  +0.97%  "

# Row 28
$ws.Range("D28").Value = "'9.79"
$ws.Range("E28").Value = "  +0.68%  "

# Row 29
$ws.Range("E29").Value = "  +0.26%  "

# Row 30
$ws.Range("D30").Value = "'0.999"
$ws.Range("E30").Value = "  -0.02%  "

# Row 31
$ws.Range("D31").Value = "'8.31"
$ws.Range("E31").Value = "  +6.08%  "

# Row 32
$ws.Range("E32").Value = "  +0.76%  "

# Row 33
$ws.Range("E33").Value = "  +5.49%  "

# Row 34
$ws.Range("D34").Value = "'495.26"
$ws.Range("E34").Value = "  -2.13%  "

# Row 35
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  +0.03%  "

# Row 36
$ws.Range("D36").Value = "'20.69"
$ws.Range("E36").Value = "  +1.58%  "

# Row 37
$ws.Range("D37").Value = "'162.83"
$ws.Range("E37").Value = "  -0.66%  "

# Row 38
$ws.Range("B38").Value = "WhiteBITCoin"
$ws.Range("C38").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D38").Value = "'20.05"
$ws.Range("E38").Value = "  +2.02%  "

# Row 39
$ws.Range("B39").Value = "PolygonEcosystemToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D39").Value = "'0.384"
$ws.Range("E39").Value = "  +3.06%  "

# Row 40
$ws.Range("D40").Value = "'0.117"
$ws.Range("E40").Value = "  +2.95%  "

# Row 41
$ws.Range("D41").Value = "'191.35"
$ws.Range("E41").Value = "  +4.98%  "

# Row 42
$ws.Range("D42").Value = "'0.105"
$ws.Range("E42").Value = "  -5.99%  "

# Row 43
$ws.Range("E43").Value = "  -0.01%  "

# Row 44
$ws.Range("D44").Value = "'0.810"
$ws.Range("E44").Value = "  +22.04%  "

# Row 45
$ws.Range("D45").Value = "'5.15"
$ws.Range("E45").Value = "  +2.44%  "

# Row 46
$ws.Range("E46").Value = "  +5.37%  "

# Row 47
$ws.Range("D47").Value = "'41.99"
$ws.Range("E47").Value = "  +4.41%  "

# Row 48
$ws.Range("E48").Value = "  -1.20%  "

# Row 49
$ws.Range("E49").Value = "  +3.90%  "

# Row 50
$ws.Range("D50").Value = "'0.608"
$ws.Range("E50").Value = "  +3.79%  "

# Row 51
$ws.Range("D51").Value = "'3.91"
$ws.Range("E51").Value = "  +3.60%  "
